$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Exams Final Score" / "SIS User ID" headers in columns A and B were
# swapped, and column A's values (Student_NN strings) are replaced with the
# numeric SIS User IDs they represent. A new "Projects Final Score" column
# (C) is also introduced with a bold header and no data yet.

$ws.Range("A1").Value = "SIS User ID"
$ws.Range("B1").Value = "Exams Final Score"

$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15

$ws.Range("C1").Value = "Projects Final Score"
$ws.Range("C1").Font.Bold = $true

$ws.Range("K13").Select()
